$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Step 1: the current "总计" sheet (sheetId=4, 4th sheet) becomes the
# new "2022-Q1" sheet. We repurpose its content in place so it keeps
# sheetId=4.
# -----------------------------------------------------------------
$quarterSheet = $wb.Worksheets.Item(4)

# Drop the old row 4 (2021-Q2 total row); we only need 3 rows (header + 2 funds)
$quarterSheet.Rows.Item(4).Delete()

# Extend the existing header style (currently on B1:D1, style index "2")
# into the new header cells E1:H1
$quarterSheet.Range("B1").Copy()
$quarterSheet.Range("E1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 1 - headers
$quarterSheet.Range("B1").Value = "基金代码"
$quarterSheet.Range("C1").Value = "基金名称"
$quarterSheet.Range("D1").Value = "基金规模"
$quarterSheet.Range("E1").Value = "股票总仓位"
$quarterSheet.Range("F1").Value = "仓位占比"
$quarterSheet.Range("G1").Value = "持有市值(亿元)"
$quarterSheet.Range("H1").Value = "仓位排名"

# Data cells in columns B:G are stored as text (not auto-converted numbers);
# column H (仓位排名) stays numeric.
$quarterSheet.Range("B2:G3").NumberFormat = "@"

# Row 2 (column A already holds 0 with the correct style)
$quarterSheet.Range("B2").Value = "610002"
$quarterSheet.Range("C2").Value = "信达澳银精华配置混合"
$quarterSheet.Range("D2").Value = "18.32"
$quarterSheet.Range("E2").Value = "71.31"
$quarterSheet.Range("F2").Value = "1.36"
$quarterSheet.Range("G2").Value = "0.2492"
$quarterSheet.Range("H2").Value = 9

# Row 3 (column A already holds 1 with the correct style)
$quarterSheet.Range("B3").Value = "610001"
$quarterSheet.Range("C3").Value = "信达澳银领先增长混合"
$quarterSheet.Range("D3").Value = "8.99"
$quarterSheet.Range("E3").Value = "94.02"
$quarterSheet.Range("F3").Value = "2.74"
$quarterSheet.Range("G3").Value = "0.2463"
$quarterSheet.Range("H3").Value = 10

$quarterSheet.Name = "2022-Q1"

# -----------------------------------------------------------------
# Step 2: add a brand-new "总计" sheet (gets the next sheetId, 5) right
# after "2022-Q1", holding the former totals plus the new 2022-Q1 row.
# -----------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("2022-Q1")
$totalSheet = $wb.Worksheets.Add($null, $afterSheet)
$totalSheet.Name = "TempTotalSheetName"

# Copy the header/index cell style (already style index "2") from the
# 2022-Q1 sheet onto the header row and column-A cells of the new sheet
$srcSheet = $wb.Worksheets.Item("2022-Q1")
$srcSheet.Range("B1").Copy()
$totalSheet.Range("B1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$srcSheet.Range("A2").Copy()
$totalSheet.Range("A2:A5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header row
$totalSheet.Range("B1").Value = "日期"
$totalSheet.Range("C1").Value = "持有数量(只)"
$totalSheet.Range("D1").Value = "持有市值(亿元)"

# Row 2 - new 2022-Q1 entry
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.5

# Row 3 - 2021-Q4 (shifted down from former row 2)
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 4
$totalSheet.Range("D3").Value = 1.01

# Row 4 - 2021-Q3 (shifted down from former row 3)
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q3"
$totalSheet.Range("C4").Value = 1
$totalSheet.Range("D4").Value = 0.18

# Row 5 - 2021-Q2 (shifted down from former row 4)
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2021-Q2"
$totalSheet.Range("C5").Value = 1
$totalSheet.Range("D5").Value = 0.15

$totalSheet.Name = "总计"
